# Update "想去人数" (want-to-go count) figures in column F across all four
# sheets to match the newly scraped data (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 540
$ws.Range("F6").Value = 2390
$ws.Range("F7").Value = 70
$ws.Range("F9").Value = 43
$ws.Range("F10").Value = 521
$ws.Range("F11").Value = 1503
$ws.Range("F13").Value = 596
$ws.Range("F14").Value = 638
$ws.Range("F15").Value = 1076
$ws.Range("F16").Value = 481
$ws.Range("F17").Value = 3383
$ws.Range("F18").Value = 405
$ws.Range("F20").Value = 3249
$ws.Range("F21").Value = 718
$ws.Range("F22").Value = 599
$ws.Range("F23").Value = 14
$ws.Range("F24").Value = 271
$ws.Range("F26").Value = 1090
$ws.Range("F28").Value = 44
$ws.Range("F29").Value = 912
$ws.Range("F30").Value = 890

# --- 演出 ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 24
$ws.Range("F8").Value = 79
$ws.Range("F13").Value = 93
$ws.Range("F19").Value = 223
$ws.Range("F20").Value = 163
$ws.Range("F21").Value = 460

# --- 本地生活 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2947
$ws.Range("F5").Value = 259
$ws.Range("F6").Value = 464

# --- 全部类型 -----------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 540
$ws.Range("F7").Value = 259
$ws.Range("F10").Value = 24
$ws.Range("F12").Value = 464
$ws.Range("F13").Value = 2390
$ws.Range("F14").Value = 70
$ws.Range("F16").Value = 43
$ws.Range("F19").Value = 79
$ws.Range("F21").Value = 521
$ws.Range("F24").Value = 1503
$ws.Range("F25").Value = 1503
$ws.Range("F28").Value = 639
$ws.Range("F29").Value = 93
$ws.Range("F31").Value = 1077
$ws.Range("F32").Value = 481
$ws.Range("F34").Value = 3383
$ws.Range("F36").Value = 3249
$ws.Range("F37").Value = 718
$ws.Range("F39").Value = 599
$ws.Range("F40").Value = 271
$ws.Range("F41").Value = 1090
$ws.Range("F44").Value = 223
$ws.Range("F45").Value = 163
$ws.Range("F46").Value = 460
$ws.Range("F48").Value = 44
$ws.Range("F49").Value = 912
$ws.Range("F50").Value = 890
